$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (2..74) down by one row (to 3..75), working from the
# bottom up so we never overwrite a row before it has been read.
for ($r = 74; $r -ge 2; $r--) {
    $src = $ws.Range("A" + $r + ":R" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":R" + ($r + 1))
    $dst.Value2 = $src.Value2
    $ws.Cells.Item($r + 1, 4).NumberFormat = $ws.Cells.Item($r, 4).NumberFormat
}

# Write the new weekly record into the now-vacated row 2.
$ws.Cells.Item(2, 1).Value2 = 11
$ws.Cells.Item(2, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(2, 3).Value2 = "Bíobío"
$ws.Cells.Item(2, 4).Value2 = 44699
$ws.Cells.Item(2, 5).Value2 = 8
$ws.Cells.Item(2, 6).Value2 = 100112012
$ws.Cells.Item(2, 7).Value2 = "Espinaca"
$ws.Cells.Item(2, 8).Value2 = "Sin especificar"
$ws.Cells.Item(2, 9).Value2 = "Primera"
$ws.Cells.Item(2, 10).Value2 = 50
$ws.Cells.Item(2, 11).Value2 = 5500
$ws.Cells.Item(2, 12).Value2 = 6000
$ws.Cells.Item(2, 13).Value2 = 5700
$ws.Cells.Item(2, 14).Value2 = "$/cuna 10 kilos"
$ws.Cells.Item(2, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(2, 16).Value2 = 570
$ws.Cells.Item(2, 17).Value2 = 10
$ws.Cells.Item(2, 18).Value2 = "Hortaliza"
